# Auto-generated Excel COM-interop script
# Applies scheduled market-data / profit-column refresh to each Leve table sheet
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2499.3
$ws.Range("J40").Value = 1998.6
$ws.Range("L40").Value = 1998.6
$ws.Range("N40").Value = -2348.6
$ws.Range("H93").Value = 48601
$ws.Range("J93").Value = 48601
$ws.Range("L93").Value = 48601
$ws.Range("N93").Value = -53593
$ws.Range("H100").Value = 1885.2858
$ws.Range("I100").Value = 1885.2858
$ws.Range("K100").Value = 1885.2858
$ws.Range("M100").Value = -1344.2858
$ws.Range("H132").Value = 1822.2927
$ws.Range("I132").Value = 1548.8387
$ws.Range("K132").Value = 4646.5161
$ws.Range("M132").Value = -2116.5161
$ws.Range("H141").Value = 696.65
$ws.Range("I141").Value = 696.65
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 2089.95
$ws.Range("L141").Value = 0
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = 3090.05

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1802.2
$ws.Range("I2").Value = 1122.5
$ws.Range("K2").Value = 1122.5
$ws.Range("M2").Value = -1009.5
$ws.Range("H32").Value = 3531.5715
$ws.Range("I32").Value = 2557.261
$ws.Range("K32").Value = 2557.261
$ws.Range("M32").Value = -2270.261
$ws.Range("H45").Value = 4455.5454
$ws.Range("I45").Value = 1752.75
$ws.Range("K45").Value = 1752.75
$ws.Range("M45").Value = -1375.75
$ws.Range("H61").Value = 1000
$ws.Range("I61").Value = 1000
$ws.Range("K61").Value = 1000
$ws.Range("M61").Value = -788
$ws.Range("H74").Value = 1001
$ws.Range("I74").Value = 1001
$ws.Range("K74").Value = 1001
$ws.Range("M74").Value = -127
$ws.Range("H77").Value = 1001
$ws.Range("I77").Value = 1001
$ws.Range("K77").Value = 5005
$ws.Range("M77").Value = -637
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").ClearContents()
$ws.Range("N95").Value = 0
$ws.Range("H96").Value = 50000
$ws.Range("J96").Value = 50000
$ws.Range("L96").Value = 50000
$ws.Range("N96").Value = -55492
$ws.Range("H116").Value = 1802.2
$ws.Range("I116").Value = 1122.5
$ws.Range("K116").Value = 1122.5
$ws.Range("M116").Value = 1171.5
$ws.Range("H122").Value = 3364.6667
$ws.Range("I122").Value = 2203.6667
$ws.Range("J122").Value = 4525.6665
$ws.Range("K122").Value = 6611.000100000001
$ws.Range("L122").Value = 13576.9995
$ws.Range("M122").Value = -4161.000100000001
$ws.Range("N122").Value = -18476.9995
$ws.Range("H132").Value = 1407.5
$ws.Range("I132").Value = 841.5
$ws.Range("K132").Value = 2524.5
$ws.Range("M132").Value = 5.5
$ws.Range("H136").Value = 1000
$ws.Range("I136").Value = 1000
$ws.Range("K136").Value = 3000
$ws.Range("M136").Value = -450

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 30000
$ws.Range("J2").Value = 30000
$ws.Range("L2").Value = 30000
$ws.Range("N2").Value = -30226
$ws.Range("H3").Value = 1802.2
$ws.Range("I3").Value = 1122.5
$ws.Range("K3").Value = 1122.5
$ws.Range("M3").Value = -1008.5
$ws.Range("H86").Value = 4917.636
$ws.Range("I86").Value = 5721.8
$ws.Range("J86").Value = 4247.5
$ws.Range("K86").Value = 5721.8
$ws.Range("L86").Value = 4247.5
$ws.Range("M86").Value = -4598.8
$ws.Range("N86").Value = -6493.5
$ws.Range("H89").Value = 4917.636
$ws.Range("I89").Value = 5721.8
$ws.Range("J89").Value = 4247.5
$ws.Range("K89").Value = 28609
$ws.Range("L89").Value = 21237.5
$ws.Range("M89").Value = -22993
$ws.Range("N89").Value = -32469.5
$ws.Range("H94").Value = 3636.7144
$ws.Range("I94").Value = 2302.1428
$ws.Range("J94").Value = 4971.2856
$ws.Range("K94").Value = 2302.1428
$ws.Range("L94").Value = 4971.2856
$ws.Range("M94").Value = -1851.1428
$ws.Range("N94").Value = -5873.2856
$ws.Range("H134").Value = 1905.5
$ws.Range("I134").Value = 1947
$ws.Range("J134").Value = 1200
$ws.Range("K134").Value = 5841
$ws.Range("L134").Value = 3600
$ws.Range("M134").Value = -3306
$ws.Range("N134").Value = -8670

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").ClearContents()
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = 0
$ws.Range("H132").Value = 4123.077
$ws.Range("I132").Value = 4225
$ws.Range("J132").Value = 2900
$ws.Range("K132").Value = 12675
$ws.Range("L132").Value = 8700
$ws.Range("M132").Value = -10145
$ws.Range("N132").Value = -13760

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 221.08333
$ws.Range("I7").Value = 22.8
$ws.Range("K7").Value = 68.40000000000001
$ws.Range("M7").Value = 43.59999999999999
$ws.Range("H56").Value = 9995
$ws.Range("I56").Value = 9995
$ws.Range("K56").Value = 9995
$ws.Range("M56").Value = -9465
$ws.Range("H113").Value = 823.7143
$ws.Range("I113").Value = 755.3333
$ws.Range("J113").Value = 875
$ws.Range("K113").Value = 2265.9999
$ws.Range("L113").Value = 2625
$ws.Range("M113").Value = -95.9998999999998
$ws.Range("N113").Value = -6965
$ws.Range("H122").Value = 728.4
$ws.Range("J122").Value = 964.3333
$ws.Range("L122").Value = 8678.9997
$ws.Range("N122").Value = -13578.9997
$ws.Range("H137").Value = 3951.5
$ws.Range("I137").Value = 1330
$ws.Range("K137").Value = 3990
$ws.Range("M137").Value = 1110

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10119.8
$ws.Range("I70").Value = 10333
$ws.Range("K70").Value = 10333
$ws.Range("M70").Value = -10063
$ws.Range("H73").Value = 10119.8
$ws.Range("I73").Value = 10333
$ws.Range("K73").Value = 10333
$ws.Range("M73").Value = -9397
$ws.Range("H132").Value = 1674.7941
$ws.Range("I132").Value = 1498.4839
$ws.Range("J132").Value = 3496.6667
$ws.Range("K132").Value = 4495.4517
$ws.Range("L132").Value = 10490.0001
$ws.Range("M132").Value = -1965.4517
$ws.Range("N132").Value = -15550.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 379.4
$ws.Range("J22").Value = 349.25
$ws.Range("L22").Value = 349.25
$ws.Range("N22").Value = -939.25
$ws.Range("H27").Value = 379.4
$ws.Range("J27").Value = 349.25
$ws.Range("L27").Value = 349.25
$ws.Range("N27").Value = -563.25
$ws.Range("H31").Value = 3515
$ws.Range("I31").Value = 3515
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 3515
$ws.Range("L31").Value = 0
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -3267
$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").ClearContents()
$ws.Range("N32").Value = 0

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 22499.5
$ws.Range("J75").Value = 22499.5
$ws.Range("L75").Value = 22499.5
$ws.Range("N75").Value = -24371.5
$ws.Range("H78").Value = 22499.5
$ws.Range("J78").Value = 22499.5
$ws.Range("L78").Value = 67498.5
$ws.Range("N78").Value = -76858.5
$ws.Range("H81").Value = 3510.3333
$ws.Range("I81").Value = 3510.3333
$ws.Range("K81").Value = 7020.6666
$ws.Range("M81").Value = -5959.6666
$ws.Range("H84").Value = 3510.3333
$ws.Range("I84").Value = 3510.3333
$ws.Range("K84").Value = 35103.333
$ws.Range("M84").Value = -29799.333
$ws.Range("H110").Value = 75644
$ws.Range("J110").Value = 75644
$ws.Range("L110").Value = 75644
$ws.Range("N110").Value = -83824
$ws.Range("H132").Value = 1308.7826
$ws.Range("I132").Value = 1308.7826
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 3926.3478
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -1396.3478
$ws.Range("H136").Value = 1666.4445
$ws.Range("I136").Value = 1666.4445
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 4999.333500000001
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -2449.333500000001
